$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-6: update values per the refreshed IFRS report (error fix)
$updates = @{
    "D2" = 2293
    "E2" = 7
    "F2" = 7
    "G2" = -54
    "H2" = -77
    "I2" = -59
    "J2" = -18
    "K2" = 2086
    "L2" = 1994
    "M2" = 92
    "N2" = 107
    "O2" = -15
    "P2" = 90
    "Q2" = 31
    "R2" = -85
    "S2" = 157
    "T2" = 82
    "U2" = -51
    "V2" = 1025
    "W2" = 0.28
    "X2" = -3.35
    "Y2" = -38.48
    "Z2" = -3.72
    "AA2" = 2161.97
    "AB2" = -55.4
    "AC2" = -327
    "AD2" = -4.29
    "AE2" = 759
    "AF2" = 1.84
    "AG2" = 20
    "AH2" = 1.43
    "AI2" = -2.2
    "AJ2" = 18000000
    "D3" = 2762
    "E3" = 60
    "F3" = 60
    "G3" = 75
    "H3" = -43
    "I3" = -8
    "J3" = -35
    "K3" = 1977
    "L3" = 1758
    "M3" = 219
    "N3" = 254
    "O3" = -34
    "P3" = 90
    "Q3" = 46
    "R3" = 483
    "S3" = -391
    "T3" = 104
    "U3" = -58
    "V3" = 593
    "W3" = 2.19
    "X3" = -1.54
    "Y3" = -4.22
    "Z3" = -2.1
    "AA3" = 801.75
    "AB3" = -59.75
    "AC3" = -42
    "AD3" = -43.39
    "AE3" = 1796
    "AF3" = 1.02
    "AG3" = 60
    "AH3" = 3.27
    "AI3" = -111.34
    "AJ3" = 18000000
    "D4" = 2938
    "E4" = 105
    "F4" = 105
    "G4" = 191
    "H4" = 152
    "I4" = 149
    "J4" = 3
    "K4" = 2006
    "L4" = 1644
    "M4" = 361
    "N4" = 358
    "O4" = 3
    "P4" = 90
    "Q4" = -19
    "R4" = -186
    "S4" = 44
    "T4" = 132
    "U4" = -151
    "V4" = 776
    "W4" = 3.56
    "X4" = 5.18
    "Y4" = 48.83
    "Z4" = 7.65
    "AA4" = 455.15
    "AB4" = 103.1
    "AC4" = 830
    "AD4" = 2.69
    "AE4" = 1991
    "AF4" = 1.12
    "AG4" = 55
    "AH4" = 2.47
    "AI4" = 6.63
    "AJ4" = 18000000
    "D5" = 3127
    "E5" = 97
    "F5" = 97
    "G5" = 45
    "H5" = 31
    "I5" = 32
    "J5" = 0
    "K5" = 2071
    "L5" = 1690
    "M5" = 381
    "N5" = 378
    "O5" = 3
    "P5" = 90
    "Q5" = 141
    "R5" = -98
    "S5" = -16
    "T5" = 113
    "U5" = 28
    "V5" = 795
    "W5" = 3.1
    "X5" = 1
    "Y5" = 8.6
    "Z5" = 1.54
    "AA5" = 443.61
    "AB5" = 116.18
    "AC5" = 176
    "AD5" = 9.720000000000001
    "AE5" = 2100
    "AF5" = 0.8100000000000001
    "AG5" = 50
    "AH5" = 2.92
    "AI5" = 28.42
    "AJ5" = 18000000
    "D6" = 3728
    "E6" = 112
    "F6" = 112
    "G6" = 84
    "H6" = 48
    "I6" = 50
    "K6" = 2124
    "L6" = 1721
    "M6" = 403
    "N6" = 401
    "P6" = 90
    "Q6" = 152
    "R6" = -210
    "S6" = 26
    "T6" = 43
    "U6" = 109
    "V6" = 827
    "W6" = 3.01
    "X6" = 1.28
    "Y6" = 12.93
    "Z6" = 2.28
    "AA6" = 427.59
    "AB6" = 154.21
    "AC6" = 280
    "AD6" = 8.52
    "AE6" = 2229
    "AF6" = 1.07
    "AG6" = 50
    "AH6" = 2.1
    "AI6" = 17.87
    "AJ6" = 18000000
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Rows 7-9: clear out now-obsolete data (only A/B/C identity columns remain)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
